# Natmi following Dr Hou advice
# Rebuild the LR-pair sheet with the 4x4 Sending-cluster x Target-cluster
# combinations (ECs, FAPs, Neutro, sCs) for Bmp6 -> Bmpr1b, replacing the
# previous 2x2x2-ish 8-row table with the full 16-row table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# shared-string lookup for the cluster names used in columns A/D
$names = @{
  20 = "ECs"
  21 = "FAPs"
  22 = "Neutro"
  23 = "sCs"
  24 = "Bmp6"
  25 = "Bmpr1b"
}

# Columns: A,D store cluster names (string); B is always "Bmp6", C is always "Bmpr1b".
# Each inner array: SendingIdx, TargetIdx, E..T (numeric columns)
$data = @(
  @(20, 20, 2, 1, 6.715603, 13.431206, 0.2130751263197567, 0.178861592654224, 1, 0.5, 0.050199, 0.100398, 0.03278033019113316, 0.02558313640211376, 0.337116554997, 1.348466219988, 0.006984672996279032, 0.004575840521972322),
  @(20, 21, 2, 1, 6.715603, 13.431206, 0.2130751263197567, 0.178861592654224, 3, 1, 0.7708836666666666, 2.312651, 0.5033929188287568, 0.5893032329676366, 5.176948664517666, 31.06169198710599, 0.1072605097679084, 0.1054037148048747),
  @(20, 22, 2, 1, 6.715603, 13.431206, 0.2130751263197567, 0.178861592654224, 1, 0.3333333333333333, 0.09074700000000001, 0.272241, 0.05925848371192178, 0.06937168705798773, 0.6094208254410001, 3.656524952646, 0.01262650890243498, 0.0124079304323021),
  @(20, 23, 2, 1, 6.715603, 13.431206, 0.2130751263197567, 0.178861592654224, 2, 1, 0.6195459999999999, 1.239092, 0.4045682672681883, 0.3157419435722618, 4.160624976237999, 16.642499904952, 0.08620343465313431, 0.05647410689507487),
  @(21, 20, 3, 1, 2.040110666666667, 6.120331999999999, 0.06472938290221787, 0.08150365120545483, 1, 0.5, 0.050199, 0.100398, 0.03278033019113316, 0.02558313640211376, 0.102411515356, 0.614469092136, 0.002121850544602991, 0.002085119026059454),
  @(21, 21, 3, 1, 2.040110666666667, 6.120331999999999, 0.06472938290221787, 0.08150365120545483, 3, 1, 0.7708836666666666, 2.312651, 0.5033929188287568, 0.5893032329676366, 1.572687991125778, 14.154191920132, 0.03258431299313168, 0.04803036515404115),
  @(21, 22, 3, 1, 2.040110666666667, 6.120331999999999, 0.06472938290221787, 0.08150365120545483, 1, 0.3333333333333333, 0.09074700000000001, 0.272241, 0.05925848371192178, 0.06937168705798773, 0.185133922668, 1.666205304012, 0.003835765082393826, 0.005654045785508197),
  @(21, 23, 3, 1, 2.040110666666667, 6.120331999999999, 0.06472938290221787, 0.08150365120545483, 2, 1, 0.6195459999999999, 1.239092, 0.4045682672681883, 0.3157419435722618, 1.263942403090667, 7.583654418543999, 0.02618745428208938, 0.02573412123984602),
  @(22, 20, 3, 1, 10.017551, 30.052653, 0.3178405490526472, 0.4002072024704813, 1, 0.5, 0.050199, 0.100398, 0.03278033019113316, 0.02558313640211376, 0.502871042649, 3.017226255894, 0.01041891814607683, 0.01023855544991068),
  @(22, 21, 3, 1, 10.017551, 30.052653, 0.3178405490526472, 0.4002072024704813, 3, 1, 0.7708836666666666, 2.312651, 0.5033929188287568, 0.5893032329676366, 7.722366445900333, 69.501298013103, 0.1599986817097467, 0.2358433982727882),
  @(22, 22, 3, 1, 10.017551, 30.052653, 0.3178405490526472, 0.4002072024704813, 1, 0.3333333333333333, 0.09074700000000001, 0.272241, 0.05925848371192178, 0.06937168705798773, 0.909062700597, 8.181564305373, 0.01883474899902457, 0.02776304880813496),
  @(22, 23, 3, 1, 10.017551, 30.052653, 0.3178405490526472, 0.4002072024704813, 2, 1, 0.6195459999999999, 1.239092, 0.4045682672681883, 0.3157419435722618, 6.206333651845998, 37.238001911076, 0.1285882001977991, 0.1263621999396475),
  @(23, 20, 2, 1, 12.7442715, 25.488543, 0.4043549417253782, 0.3394275536698397, 1, 0.5, 0.050199, 0.100398, 0.03278033019113316, 0.02558313640211376, 0.6397496850285, 2.558998740114, 0.01325488850417431, 0.008683621404171298),
  @(23, 21, 2, 1, 12.7442715, 25.488543, 0.4043549417253782, 0.3394275536698397, 3, 1, 0.7708836666666666, 2.312651, 0.5033929188287568, 0.5893032329676366, 9.824350742915499, 58.94610445749299, 0.20354941435797, 0.2000257547359326),
  @(23, 22, 2, 1, 12.7442715, 25.488543, 0.4043549417253782, 0.3394275536698397, 1, 0.3333333333333333, 0.09074700000000001, 0.272241, 0.05925848371192178, 0.06937168705798773, 1.1565044058105, 6.939026434863, 0.0239614607280684, 0.02354666203204245),
  @(23, 23, 2, 1, 12.7442715, 25.488543, 0.4043549417253782, 0.3394275536698397, 2, 1, 0.6195459999999999, 1.239092, 0.4045682672681883, 0.3157419435722618, 7.895662430738999, 31.582649722956, 0.1635891781351655, 0.1071715154976934)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $names[[int]$row[0]]
    $ws.Cells.Item($r, 2).Value = "Bmp6"
    $ws.Cells.Item($r, 3).Value = "Bmpr1b"
    $ws.Cells.Item($r, 4).Value = $names[[int]$row[1]]

    for ($j = 2; $j -lt $row.Length; $j++) {
        $col = $j + 3   # E=5 .. T=20
        $ws.Cells.Item($r, $col).Value = $row[$j]
    }
}
